$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.2169117647058824
$ws.Range("C2").Value = 0.5183823529411765
$ws.Range("J2").Value = 0.01102941176470588
$ws.Range("P2").Value = 0.1470588235294118
$ws.Range("S2").Value = 0.1066176470588235
$ws.Range("B3").Value = 0.0131578947368421
$ws.Range("C3").Value = 0.04605263157894737
$ws.Range("J3").Value = 0.02631578947368421
$ws.Range("P3").Value = 0.7105263157894737
$ws.Range("S3").Value = 0.2039473684210526
$ws.Range("P4").Value = 0.5208333333333334
$ws.Range("S4").Value = 0.4791666666666667
$ws.Range("B6").Value = 0.0611353711790393
$ws.Range("D6").Value = 0.01310043668122271
$ws.Range("E6").Value = 0.004366812227074236
$ws.Range("F6").Value = 0.04366812227074236
$ws.Range("J6").Value = 0.2008733624454148
$ws.Range("O6").Value = 0.02620087336244541
$ws.Range("Q6").Value = 0.1615720524017467
$ws.Range("R6").Value = 0.07423580786026202
$ws.Range("S6").Value = 0.4148471615720524
$ws.Range("B7").Value = 0.07471264367816093
$ws.Range("D7").Value = 0.01724137931034483
$ws.Range("F7").Value = 0.06321839080459771
$ws.Range("J7").Value = 0.132183908045977
$ws.Range("O7").Value = 0.01724137931034483
$ws.Range("Q7").Value = 0.1839080459770115
$ws.Range("R7").Value = 0.1379310344827586
$ws.Range("S7").Value = 0.3735632183908046
$ws.Range("B8").Value = 0.08997955010224949
$ws.Range("D8").Value = 0.016359918200409
$ws.Range("E8").Value = 0.002044989775051125
$ws.Range("F8").Value = 0.05930470347648262
$ws.Range("J8").Value = 0.08793456032719836
$ws.Range("O8").Value = 0.0245398773006135
$ws.Range("Q8").Value = 0.2004089979550102
$ws.Range("R8").Value = 0.07975460122699386
$ws.Range("S8").Value = 0.4396728016359918
$ws.Range("B9").Value = 0.1141552511415525
$ws.Range("D9").Value = 0.0182648401826484
$ws.Range("E9").Value = 0.0045662100456621
$ws.Range("F9").Value = 0.0730593607305936
$ws.Range("J9").Value = 0.1050228310502283
$ws.Range("O9").Value = 0.0273972602739726
$ws.Range("Q9").Value = 0.1780821917808219
$ws.Range("R9").Value = 0.0730593607305936
$ws.Range("S9").Value = 0.4063926940639269
$ws.Range("B10").Value = 0.0977891156462585
$ws.Range("D10").Value = 0.02465986394557823
$ws.Range("E10").Value = 0.001700680272108843
$ws.Range("F10").Value = 0.08248299319727891
$ws.Range("J10").Value = 0.08928571428571429
$ws.Range("O10").Value = 0.02380952380952381
$ws.Range("Q10").Value = 0.2091836734693878
$ws.Range("R10").Value = 0.0858843537414966
$ws.Range("S10").Value = 0.3852040816326531
$ws.Range("G11").Value = 0.1879699248120301
$ws.Range("J11").Value = 0.09398496240601503
$ws.Range("K11").Value = 0.2368421052631579
$ws.Range("L11").Value = 0.4699248120300752
$ws.Range("S11").Value = 0.0112781954887218
$ws.Range("G12").Value = 0.7244094488188977
$ws.Range("J12").Value = 0.2440944881889764
$ws.Range("K12").Value = 0.007874015748031496
$ws.Range("L12").Value = 0.01574803149606299
$ws.Range("S12").Value = 0.007874015748031496
$ws.Range("G13").Value = 0.6666666666666666
$ws.Range("J13").Value = 0.2745098039215687
$ws.Range("S13").Value = 0.05882352941176471
$ws.Range("F15").Value = 0.01351351351351351
$ws.Range("H15").Value = 0.1621621621621622
$ws.Range("I15").Value = 0.04954954954954955
$ws.Range("J15").Value = 0.3423423423423423
$ws.Range("K15").Value = 0.07657657657657657
$ws.Range("M15").Value = 0.01351351351351351
$ws.Range("O15").Value = 0.03603603603603604
$ws.Range("S15").Value = 0.3063063063063063
$ws.Range("F16").Value = 0.01796407185628742
$ws.Range("H16").Value = 0.1856287425149701
$ws.Range("I16").Value = 0.1197604790419162
$ws.Range("J16").Value = 0.3712574850299401
$ws.Range("K16").Value = 0.09580838323353294
$ws.Range("M16").Value = 0.01197604790419162
$ws.Range("N16").Value = 0.005988023952095809
$ws.Range("O16").Value = 0.07784431137724551
$ws.Range("S16").Value = 0.1137724550898204
$ws.Range("F17").Value = 0.02654867256637168
$ws.Range("H17").Value = 0.1880530973451327
$ws.Range("I17").Value = 0.1172566371681416
$ws.Range("J17").Value = 0.3849557522123894
$ws.Range("K17").Value = 0.05973451327433629
$ws.Range("M17").Value = 0.01769911504424779
$ws.Range("N17").Value = 0.00663716814159292
$ws.Range("O17").Value = 0.05973451327433629
$ws.Range("S17").Value = 0.1393805309734513
$ws.Range("F18").Value = 0.01015228426395939
$ws.Range("H18").Value = 0.1928934010152284
$ws.Range("I18").Value = 0.1370558375634518
$ws.Range("J18").Value = 0.3807106598984771
$ws.Range("K18").Value = 0.08121827411167512
$ws.Range("M18").Value = 0.02030456852791878
$ws.Range("O18").Value = 0.07614213197969544
$ws.Range("S18").Value = 0.1015228426395939
$ws.Range("F19").Value = 0.01278195488721805
$ws.Range("H19").Value = 0.2270676691729323
$ws.Range("I19").Value = 0.08195488721804511
$ws.Range("J19").Value = 0.3699248120300752
$ws.Range("K19").Value = 0.09323308270676692
$ws.Range("M19").Value = 0.02631578947368421
$ws.Range("O19").Value = 0.05864661654135338
$ws.Range("S19").Value = 0.1300751879699248